$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the available-quantity (Stock dispo) column with the new reduced values
$ws.Range("B2").Value = 6
$ws.Range("B3").Value = 40
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 11
$ws.Range("B6").Value = 140

# Move the active selection to D5, matching the final cursor position
$ws.Range("D5").Select()
